$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Container With Most Water
$s34 = @"
Container with most water
"@
$ws.Range("B12").Value = $s34
$ws.Range("C12").Value = "Two Pointers"
$ws.Range("D12").Value = "List"
$s35 = @"
1. Initialize two pointers at the two ends of the array and calculate the area formed by the lines they point to
2. Move the pointer pointing to the shorter line inward to potentially find a taller line and a larger area
3. The two pointer approach makes sure an O(n) time complexity by making a single pass through the array
"@
$e12 = $ws.Range("E12")
$e12.Value = $s35
$e12.WrapText = $true
$ws.Rows.Item(12).RowHeight = 85

# Row 9: Valid Anagram - add notes
$s36 = @"
1. Counting: Traverse both strings and count character occurences 
2. Comparison: Compare the character count for both strings
3. Alternatively sort both strings and do a direct comparison
4. Counting is O(n), sorting is O(n*logn)
"@
$e9 = $ws.Range("E9")
$e9.Value = $s36
$e9.WrapText = $true
$ws.Rows.Item(9).RowHeight = 68

# Row 10: Two Sum Two - add notes
$s37 = @"
1. Pointer initialization: Initialize two pointers at the beiginning and end of the array.
2. Pointer movement: Move the left pointer inwards or the right pointer inwards based on the sum comparison with the target
3. Edge cases: handle cases where array has fewer than two elements 
4. O(n) due to a single pass with two pointers
"@
$e10 = $ws.Range("E10")
$e10.Value = $s37
$e10.WrapText = $true
$ws.Rows.Item(10).RowHeight = 85

# Row 13: Best Time to Buy and Sell Stock
$s38 = @"
Best time to buy and sell stock
"@
$ws.Range("B13").Value = $s38
$ws.Range("C13").Value = "Sliding Window"
$ws.Range("D13").Value = "List"
$s39 = @"
Note: try this problem out with greedy approach the next time you see this
1. Learn about greedy algorithm.
"@
$e13 = $ws.Range("E13")
$e13.Value = $s39
$e13.WrapText = $true
$ws.Rows.Item(13).RowHeight = 34

# Row 14: Longest Substring Without Repeating Characters
$s40 = @"
Longest substring without repeating characters
"@
$ws.Range("B14").Value = $s40
$ws.Range("C14").Value = "Sliding Window"
$s41 = @"
Hashmap or Set
"@
$ws.Range("D14").Value = $s41
$s42 = @"
1. Use sliding window technique to maintain a substring with unique chars.
2. window is defined by two pointers: start and end / left and right
3. as we traverse based on conditions we adjust the window size dynamically
4. we have the hashmap to store key ( char ) and value ( index position of that char )
5. if a repeating character is found we move the position of the start pointer next to the last seen index of the repeating char. then update the hashmap with current index of the char
6. we calculate the current window size : end - start + 1 , and then update the max_len( max len of the substring).
"@
$e14 = $ws.Range("E14")
$e14.Value = $s42
$e14.WrapText = $true
$ws.Rows.Item(14).RowHeight = 136

# View state: zoom + selection
$excel.ActiveWindow.Zoom = 203
[void]$ws.Range("E14").Select()
